$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Commit: "all the test file got parameterised. and created enhanced
# assertion for all test files. now assertion looks good"
#
# The Test Objective (col C) for each case got re-derived/reshuffled, every
# row's Pre-Conditions/Test Data/Test Steps/Expected Result (cols D-G) were
# parameterised down to one shared "Reset Pass" template, the PASSED actual
# result/status moved from TC_013 (row 14) to TC_009 (row 10), and a brand
# new case (TC_024, row 25) was appended at the bottom.
# ---------------------------------------------------------------------------

# New Test Objective text per row (2..25); A/B (serial no. / TC id) stay as-is.
$objectives = @{
    2  = 'Verify forgot password link is visible on login page.'
    3  = 'Verify forgot password page is visible after clicking the link.'
    4  = 'Verify email required validation is visible.'
    5  = 'Verify email input section''s Next button is disabled.'
    6  = 'Verify invalid email format validation is visible.'
    7  = 'Verify unregistered email validation is visible.'
    8  = 'Verify public domain email validation is visible.'
    9  = 'Verify non-verified email validation is visible.'
    10 = 'Verify forgot password page''s elements are visible.'
    11 = 'Verify reset password page''s elements are visible.'
    12 = 'Verify resend OTP countdown is visible.'
    13 = 'Verify valid email navigation to reset password page.'
    14 = 'Verify empty fields validation message is visible.'
    15 = 'Verify OTP input accepts only numbers.'
    16 = 'Verify max attempt limit exceeded error is visible.'
    17 = 'Verify OTP input limit validation with less than 6 digits.'
    18 = 'Verify OTP input limit validation with more than 6 digits.'
    19 = 'Verify invalid OTP validation is visible.'
    20 = 'Verify password complexity validation is visible.'
    21 = 'Verify password mismatch validation is visible.'
    22 = 'Verify show/hide password functionality works.'
    23 = 'Verify resend OTP functionality works.'
    24 = 'Verify back button functionality on reset password page.'
    25 = 'Verify back button functionality on email input page.'
}

# Shared parameterised pre-condition / test data / steps / expected-result
# block now used by every row.
$resetD = 'User is on the reset pass page'
$resetE = 'Reset Pass test data'
$resetF = "1. Navigate to reset pass page`n2. Perform required actions`n3. Verify expected behavior"
$resetG = 'Reset Pass functionality should work as expected'

# New row 25/TC_024 didn't exist before - clone row 2's cell formatting (incl.
# the 80pt custom row height) onto it, then write its serial number + id.
$ws.Range("A2:I2").Copy()
$ws.Range("A25:I25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(25).RowHeight = 80

# A25 must stay text ("24"), matching the zero-padded text siblings above it -
# a plain .Value assignment would get auto-coerced to the number 24, so force
# text storage via a quote-prefixed formula, then re-stamp the clean (no
# quote-prefix) number format/style from a sibling cell.
$ws.Range("A25").Formula = "'24"
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B25").Value = "TC_024"

foreach ($r in 2..25) {
    $ws.Range("C$r").Value = $objectives[$r]
    $ws.Range("D$r").Value = $resetD
    $ws.Range("E$r").Value = $resetE
    $ws.Range("F$r").Value = $resetF
    $ws.Range("G$r").Value = $resetG
}

# Actual Result / Test Status default back to "not executed" for every row...
foreach ($r in 2..25) {
    $ws.Range("H$r").Value = "Test not executed"
    $ws.Range("I$r").Value = "Not Run"
}

# ...then copy the original PASSED (green) formatting off row 14 before we
# overwrite its text, and move the PASSED result onto row 10 (TC_009).
$ws.Range("I14").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H10").Value = "Reset Pass functionality verified"
$ws.Range("I10").Value = "PASSED"

# Row 14 (TC_013) reverts to the plain "Not Run" formatting (copy it from a
# row that never had the PASSED styling).
$ws.Range("I2").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H14").Value = "Test not executed"
$ws.Range("I14").Value = "Not Run"
